$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.809971570968628
$ws.Range("B1").Value = 6.129008769989014
$ws.Range("C1").Value = 5.186582088470459
$ws.Range("D1").Value = 6.031337738037109
$ws.Range("E1").Value = 3.940948724746704
